# The document currently starts with:
#   P1 (style Heading1): a single run "On Pilgrimage - March 1951"
#   P2 (bold run):        a single run "By Dorothy Day"
# (P1 is wrapped in a bookmarkStart/bookmarkEnd pair named
#  "on-pilgrimage---march-1951".)
#
# Pandoc's docx writer renders a markdown title block as a "Title"
# styled paragraph and an "Authors" styled paragraph, each word (and
# each inter-word space / punctuation token) as its own run. Reproduce
# that structure here:
#   P1 (style Title):   "On" " " "Pilgrimage" " " "-" " " "March" " " "1951"
#   P2 (style Authors):  "Dorothy" " " "Day"

$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function New-RunsXml($words) {
    $runs = ""
    foreach ($tok in $words) {
        $runs += '<w:r><w:t xml:space="preserve">' + $tok + '</w:t></w:r>'
    }
    return $runs
}

$titleRuns  = New-RunsXml @("On", " ", "Pilgrimage", " ", "-", " ", "March", " ", "1951")
$authorRuns = New-RunsXml @("Dorothy", " ", "Day")

$titleXml  = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:pStyle w:val="Title"/></w:pPr>' + $titleRuns + '</w:p>'
$authorXml = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:pStyle w:val="Authors"/></w:pPr>' + $authorRuns + '</w:p>'

# Replace paragraph 1 (the title heading) in place, keeping it as
# paragraph 1 so the bookmark that used to wrap it still wraps the new
# title paragraph.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertXML($titleXml)

# Replace paragraph 2 (the "By Dorothy Day" byline) in place.
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertXML($authorXml)

Write-Host "Paragraph count:" $d.Paragraphs.Count
Write-Host "P1 style/text:" $d.Paragraphs.Item(1).Range.ParagraphStyle.NameLocal "|" $d.Paragraphs.Item(1).Range.Text
Write-Host "P2 style/text:" $d.Paragraphs.Item(2).Range.ParagraphStyle.NameLocal "|" $d.Paragraphs.Item(2).Range.Text
Write-Host "P3 text:" $d.Paragraphs.Item(3).Range.Text
